# modify SLG building config
#
# Summary of changes (per the target diff):
#  - Property sheet: add a new row 13 (LoadPropertyFinish / int / true / true /
#    true / 0 / 0 / Friend), extend dimension, fix up the F-column dropdown
#    data validation so it is a single contiguous range again, and move the
#    sheet's remembered selection to C26.
#  - Record_BuildingList sheet: bump C2 from 6 to 8, move selection to G10 and
#    make this sheet the active tab.
#  - Record_BuildingProduce sheet: move its remembered selection to F28 (it is
#    no longer the active tab once Record_BuildingList is activated below).

$wb = $excel.ActiveWorkbook

$wsProperty = $wb.Worksheets.Item("Property")
$wsBuildingList = $wb.Worksheets.Item("Record_BuildingList")
$wsBuildingProduce = $wb.Worksheets.Item("Record_BuildingProduce")

# --- Property: append row 13 ------------------------------------------------
$wsProperty.Range("A13").Value = "LoadPropertyFinish"

$wsProperty.Range("B13").NumberFormat = "@"
$wsProperty.Range("B13").Value = "int"

$wsProperty.Range("C13").Value = $true
$wsProperty.Range("D13").Value = $true
$wsProperty.Range("E13").Value = $true

$wsProperty.Range("G13").Value = 0
$wsProperty.Range("H13").Value = 0

$wsProperty.Range("I13").NumberFormat = "@"
$wsProperty.Range("I13").Value = "Friend"

# Re-merge the F-column TRUE/FALSE dropdown validation: it used to be split
# into "F2:F12" + "F13:F1048576" (an artifact of row 13 not existing yet) and
# should now cover "F2:F1048576" as a single rule.
$fRange = $wsProperty.Range("F2:F1048576")
$fRange.Validation.Delete() | Out-Null
$fRange.Validation.Add(3, 1, 1, '"TRUE,FALSE"') | Out-Null

# Move the saved selection on the Property sheet.
$wsProperty.Range("C26").Select() | Out-Null

# --- Record_BuildingList: bump SaveInterval-ish C2 6 -> 8 -------------------
$wsBuildingList.Range("C2").Value = 8

# --- Record_BuildingProduce: just relocate the saved selection -------------
$wsBuildingProduce.Range("F28").Select() | Out-Null

# --- Activate Record_BuildingList last so it becomes the active tab --------
$wsBuildingList.Activate()
$wsBuildingList.Range("G10").Select() | Out-Null
